$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Narrow the status-date columns ---
# Target stored column width (OOXML) is 13.4101845877511 characters;
# this runtime derives the stored width from ColumnWidth via
# stored = ColumnWidth + 5/6, so back it out from the target.
$targetColumnWidth = 13.4101845877511 - (5/6)

$overview.Range("E1").ColumnWidth = $targetColumnWidth
$overview.Range("F1").ColumnWidth = $targetColumnWidth

$zhcn.Range("C1").ColumnWidth = $targetColumnWidth
$dede.Range("C1").ColumnWidth = $targetColumnWidth
